# "simple model workes with 24 timesteps"
#
# - General Data: a couple of input values change, and one label is renamed.
# - Demand: a new totals ("SUM") row is added at the bottom of the table.
# - A brand-new "Tabelle1" sheet is added at the end of the workbook, holding
#   a small supply-balance summary, and becomes the active sheet.
# - Several sheets simply end up with a different cell selected (because that
#   is where the user happened to click while making the edits above).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Demand" sheet - add a SUM row under the data (row 26).
#    Do this first so the shared string "SUM" is appended to the shared
#    string table before the other new strings below.
# ---------------------------------------------------------------------------
$wsDemand = $wb.Worksheets.Item("Demand")
$wsDemand.Range("A26").Value = "SUM"
$wsDemand.Range("B26").Formula = "=SUM(B2:B25)"
$wsDemand.Range("C26:E26").Formula = "=SUM(C2:C25)"

# ---------------------------------------------------------------------------
# 2) "General Data" sheet tweaks
# ---------------------------------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("General Data")
$wsGeneral.Range("C4").Value = 20
$wsGeneral.Range("B14").Value = "η_battery_gas"
$wsGeneral.Range("C15").Value = 2

# ---------------------------------------------------------------------------
# 3) New sheet "Tabelle1" with a small supply-balance summary table.
# ---------------------------------------------------------------------------
$wsNew = $wb.Worksheets.Add()
$wsNew.Name = "Tabelle1"

$wsNew.Range("A9").Value = "Thermal supply"
$wsNew.Range("A1").Value = "Supply"
$wsNew.Range("A10").Value = "Electric supply"

$wsNew.Range("B1").Value = "Contractor"
$wsNew.Range("C1").Value = "Self financed"

$wsNew.Range("A2").Value = "PV"
$wsNew.Range("B2").Value = 8

$wsNew.Range("A3").Value = "ST"
$wsNew.Range("B3").Value = 0

$wsNew.Range("A4").Value = "HP"

$wsNew.Range("A5").Value = "Grid"
$wsNew.Range("B5").Value = 82

$wsNew.Range("A6").Value = "DH"
$wsNew.Range("B6").Value = 180

$wsNew.Range("A7").Value = "Gas"

$wsNew.Range("B9").Formula = "=B6+B3"
$wsNew.Range("B10").Formula = "=B2+B5"

$wsNew.Columns.Item(1).AutoFit() | Out-Null

# Move it to the end of the tab strip and make it the active sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew.Move($null, $lastSheet)
$wb.Worksheets.Item("Tabelle1").Activate()
$wb.Worksheets.Item("Tabelle1").Range("B6").Select()

# ---------------------------------------------------------------------------
# 4) Selection bookkeeping on the other sheets (the cell the user last
#    clicked on each sheet while working).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Sets").Range("K2:K5").Select()
$wb.Worksheets.Item("General Data").Range("C15").Select()
$wb.Worksheets.Item("Costs new investments").Range("B3").Select()
$wb.Worksheets.Item("Costs default system").Range("E4").Select()
$wb.Worksheets.Item("Demand").Range("E26").Select()
$wb.Worksheets.Item("irradiation_winter").Range("E15").Select()

# Re-activate Tabelle1 last so it ends up the active/selected tab, matching
# the recorded workbook view state.
$wb.Worksheets.Item("Tabelle1").Activate()
$wb.Worksheets.Item("Tabelle1").Range("B6").Select()
